$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

$updates = @{
    "B21"  = "Details of any changes the proposed development would make to existing access arrangements or public rights of way"
    "B26"  = "Name and contact information if an agent is being used."
    "B30"  = "Name and contact information if an agent is being used."
    "B38"  = "Telephone number and email address of the applicant."
    "B42"  = "Name and contact information for the parties making the application."
    "B48"  = "How any natural habitats on the development site will be improved by the proposed works."
    "B49"  = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
    "B50"  = "Details of any conflict of interest that may exist between the applicant and planning authority."
    "B53"  = "Signed and dated verification of the application's accuracy."
    "B56"  = "What materials are being used for the proposed development"
    "B64"  = "Who will be affected by the proposal and whether they have been notified, such as agricultural tenants"
    "B79"  = "Details of any changes the proposed development would make to parking facilities."
    "B81"  = "Details of pre-application advice received from the planning authority"
    "B86"  = "What development, works or change of use is proposed"
    "B91"  = "Where the proposed development will be built."
    "B100" = "Information to help the planning authority arrange a site visit"
    "B106" = "Details of trees and/or hedges that will be affected by the proposed development"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
